$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-01 Tuesday" "2024-10-02 Wednesday"

Replace-Text "493×3=1479" "166×4=664"
Replace-Text "209×3=627" "188×3=564"
Replace-Text "994×3=2982" "972×3=2916"
Replace-Text "945×7=6615" "236×8=1888"
Replace-Text "809×2=1618" "727×8=5816"
Replace-Text "690×9=6210" "949×6=5694"
Replace-Text "757×3=2271" "142×5=710"
Replace-Text "533×5=2665" "705×9=6345"
Replace-Text "443×2=886" "187×7=1309"
Replace-Text "684×6=4104" "420×7=2940"
Replace-Text "829×4=3316" "751×3=2253"
Replace-Text "129×2=258" "462×6=2772"
Replace-Text "817×5=4085" "907×4=3628"
Replace-Text "926×9=8334" "222×3=666"
Replace-Text "607×4=2428" "707×4=2828"
Replace-Text "614×4=2456" "273×5=1365"
Replace-Text "478×5=2390" "205×9=1845"
Replace-Text "475×2=950" "976×5=4880"
Replace-Text "519×7=3633" "540×6=3240"
Replace-Text "665×4=2660" "268×3=804"
Replace-Text "837×9=7533" "390×9=3510"
Replace-Text "842×5=4210" "524×6=3144"
Replace-Text "960×4=3840" "220×7=1540"
Replace-Text "278×2=556" "847×8=6776"
Replace-Text "952×8=7616" "729×5=3645"
